$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text representation (avoid Excel auto-converting
# numeric/date-looking strings like "69.50" or "0.0550" into numbers).
$cells = @{
    "D2" = "34.796.04"
    "E2" = "  -1.00%  "
    "D3" = "1.829.43"
    "E3" = "  +0.55%  "
    "E4" = "  +0.30%  "
    "D5" = "230.88"
    "E5" = "  -1.10%  "
    "E6" = "  +0.59%  "
    "E7" = "  +0.27%  "
    "D8" = "39.52"
    "E8" = "  -4.78%  "
    "D9" = "0.325"
    "E9" = "  +0.12%  "
    "D10" = "0.0683"
    "E10" = "  -0.63%  "
    "D11" = "0.0987"
    "E11" = "  -1.32%  "
    "D12" = "2.092.58"
    "E12" = "  +0.50%  "
    "D13" = "1.852.22"
    "E13" = "  +1.76%  "
    "D14" = "11.29"
    "E14" = "  +1.49%  "
    "E15" = "  +1.07%  "
    "D16" = "4.62"
    "E16" = "  -1.41%  "
    "D17" = "34.754.90"
    "E17" = "  -0.92%  "
    "D18" = "69.50"
    "E18" = "  -0.33%  "
    "D19" = "0.0₃0785"
    "E19" = "  -1.12%  "
    "D20" = "240.36"
    "E20" = "  +0.12%  "
    "D21" = "12.12"
    "E21" = "  +1.82%  "
    "D22" = "4.67"
    "E22" = "  -0.15%  "
    "E23" = "  +0.29%  "
    "E24" = "  -0.25%  "
    "D25" = "171.47"
    "E25" = "  -0.89%  "
    "D26" = "7.74"
    "E26" = "  -1.72%  "
    "E27" = "  +2.08%  "
    "D28" = "17.32"
    "E28" = "  -1.38%  "
    "E29" = "  -7.13%  "
    "E30" = "  +0.31%  "
    "D31" = "0.0550"
    "D32" = "3.92"
    "E32" = "  -3.48%  "
    "D33" = "3.93"
    "E33" = "  -1.61%  "
    "E34" = "  +3.03%  "
    "E35" = "  +6.82%  "
    "D36" = "1.43"
    "E36" = "  +11.82%  "
    "D37" = "0.697"
    "E37" = "  +1.51%  "
    "D38" = "91.11"
    "E38" = "  -2.45%  "
    "E39" = "  +6.01%  "
    "D40" = "1.337.88"
    "E40" = "  +2.03%  "
    "E41" = "  -0.85%  "
    "D42" = "14.55"
    "E42" = "  -1.71%  "
    "D43" = "2.42"
    "E43" = "  -1.96%  "
    "E44" = "  -3.40%  "
    "E45" = "  -0.64%  "
    "E46" = "  +1.95%  "
    "E47" = "  -1.95%  "
    "D48" = "2.007.11"
    "E48" = "  +0.51%  "
    "E49" = "  +0.32%  "
    "E50" = "  +2.83%  "
    "E51" = "  +13.23%  "
}

foreach ($ref in $cells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$ref]
}
